$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing rows (old rows 13-20), shrinking the
# used range down to A1:B12.
$ws.Range("A13:B20").EntireRow.Delete()

# Years + values shift: row 2 used to hold 2001, now holds 2010, etc.
$years = @("2010年","2011年","2013年","2014年","2015年","2016年","2017年","2018年","2019年","2020年","2021年")
$values = @(138904.7, 604171.3204, 465044.8, 154652.8717, 115737.3, 102799.5, 37946.4, 39367.2, 346185.3979, 205722.5718, 206046)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $years[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
